$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells J1:L1 - copy formatting from existing header (I1) so they
# share the same bold/border/centered style instead of getting a new style index.
$ws.Range("I1").Copy()
$ws.Range("J1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J1").Value = "Tanggal_Bayar"
$ws.Range("K1").Value = "No_Resi"
$ws.Range("L1").Value = "Status_Pengiriman"

# New data row 2 - plain (unstyled) text values.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4545454444444444"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "BG8888YY"
$ws.Range("C2").Value = "Siti Aminah"
$ws.Range("D2").Value = "Palembang"
$ws.Range("E2").Value = "06-08-2025 14:37"
$ws.Range("F2").Value = "Pos Indonesia"
$ws.Range("G2").Value = "RESI975880"
$ws.Range("H2").Value = "Diproses"
$ws.Range("I2").Value = "08-08-2025 14:37"
$ws.Range("J2").Value = "06-08-2025 14:37"
$ws.Range("K2").Value = "RESI975880"
$ws.Range("L2").Value = "Diproses"
